# Generate Report for Archive
#
# 1) The localization status "Ready for handoff" has moved on to
#    "In Translation" everywhere it appears (Overview sheet's per-language
#    status columns, plus the Status column on each per-language sheet).
# 2) The two "Status"-ish columns on the Overview sheet (zh-cn / de-de) and
#    the Status column on each language sheet are narrowed to match the new,
#    shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn sheet: Status column (C) ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# --- de-de sheet: Status column (C) ---
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Narrow the columns that held the longer "Ready for handoff" text ---
# ColumnWidth is expressed in characters; use a value that lands on the
# narrower width the report now needs.
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
